# [PHOENIX-6082] completed Forward/Close Grievance
#
# The "approvalDetails" sheet had two shared-string values that included a
# "S." prefix in front of "Ravindra Babu". Those values are updated to drop
# the "S." prefix. The active selection on the sheet is also moved from
# B22 to D18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("approvalDetails")

# Update the approver name for the "commissioner" row (D5) and the
# "commissioner1" row (D6), which also cascades to D14 since it shares the
# same text as D6.
$ws.Range("D5").Value = "Ravindra Babu ~ ADM_Commissioner_1"
$ws.Range("D6").Value = "Ravindra Babu/ADM_Commissioner_1"
$ws.Range("D14").Value = "Ravindra Babu/ADM_Commissioner_1"

# Move the selection to D18, matching the latest edit location.
$ws.Range("D18").Select()
